# Update countries & provincias Spain
# Refresh the COVID-19 stats scrape (new run at 06:16 instead of 04:59) and
# reorder Bahamas/Guinea-Bisau in the country list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 06:16"

# --- India (row 6) ------------------------------------------------------
$ws.Range("B6").Value = 3691166
$ws.Range("C6").Value = 3227
$ws.Range("D6").Value = 2839882
$ws.Range("E6").Value = 785849

# --- Honduras (row 50) ---------------------------------------------------
$ws.Range("B50").Value = 61014
$ws.Range("C50").Value = 840
$ws.Range("D50").Value = 10396
$ws.Range("E50").Value = 48745
$ws.Range("G50").Value = 15
$ws.Range("H50").Value = 1873

# --- Bahamas / Guinea-Bisau reorder (rows 138 & 139) ----------------------
# Bahamas now sorts before Guinea-Bisau; row 138 gets Bahamas' refreshed
# figures and row 139 keeps Guinea-Bisau's previous (unchanged) figures.
$ws.Range("A138").Value = "Bahamas"
$ws.Range("B138").Value = 2217
$ws.Range("C138").Value = 0
$ws.Range("D138").Value = 801
$ws.Range("E138").Value = 1366
$ws.Range("H138").Value = 50

$ws.Range("A139").Value = "Guinea-Bisau"
$ws.Range("B139").Value = 2205
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 1127
$ws.Range("E139").Value = 1044
$ws.Range("H139").Value = 34

# --- Islas Turcas y Caicos (row 173) --------------------------------------
$ws.Range("B173").Value = 508
$ws.Range("C173").Value = 1
$ws.Range("G173").Value = 1
$ws.Range("H173").Value = 4

# --- San Martin (Parte Holandesa) (row 175) -------------------------------
$ws.Range("B175").Value = 476
$ws.Range("C175").Value = 13
$ws.Range("D175").Value = 200
$ws.Range("E175").Value = 259

# --- Mongolia (row 183) ---------------------------------------------------
$ws.Range("B183").Value = 304
$ws.Range("C183").Value = 3
$ws.Range("E183").Value = 9

# --- Butan (row 186) -------------------------------------------------------
$ws.Range("B186").Value = 225
$ws.Range("C186").Value = 1
$ws.Range("D186").Value = 142
$ws.Range("E186").Value = 83
